# Currency.xlsx edit:
#  1. Column C (rows 2-201) is flattened to a constant value of 500,
#     replacing the previous formulas/values (5000, 4975, SUM(...) chains, etc.)
#  2. The sheet's active cell/selection moves from D2 to L22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite the whole Stardust-earned column with the flat value 500.
# Assigning a plain scalar to a multi-cell range fills every cell and
# clears any existing formulas, matching the target workbook state.
$ws.Range("C2:C201").Value = 500

# Update the saved selection/active cell to L22.
$ws.Range("L22").Select()
